$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.579.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.62"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.87%  "

$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.54"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("E14").Value = "  -3.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.895.78"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.482.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.447.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.68"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.35%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "641.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.574.55"
$ws.Range("D27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0951"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.78"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.74%  "

$ws.Range("E32").Value = "  -3.85%  "

$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.93"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.45"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.362"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("E42").Value = "  -2.94%  "

$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0308"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.34%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("E47").Value = "  -2.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.602"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("E51").Value = "  -1.95%  "
